# GanttChart.xlsx - "Project Planner" sheet update
# Bring every in-progress task's "Percent Complete" (column G) up to 100%,
# and move the highlighted reporting period (H2, drives the " Period
# Highlight" conditional formatting) forward from period 9 to period 44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Percent Complete column (G) -> 100% for every task that wasn't already done ---
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G11:G36").Value = 1

# --- Highlighted period selector ---
$ws.Range("H2").Value = 44

# --- Restore the view the workbook was left in: zoomed out a touch, ---
# --- scrolled down, with F35 as the active selection ---
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("F35").Select()

Write-Host "Updated Percent Complete values and selected period."
